$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (not numeric-looking) -- safe to assign directly.
$textUpdates = @(
    @('D2', '42.571.00'),
    @('E2', '  -0.85%  '),
    @('D3', '2.231.49'),
    @('E3', '  -0.41%  '),
    @('E4', '  -0.13%  '),
    @('E5', '  -1.18%  '),
    @('E6', '  +6.37%  '),
    @('E7', '  -0.48%  '),
    @('E8', '  -0.38%  '),
    @('E9', '  -1.31%  '),
    @('E10', '  -6.23%  '),
    @('E11', '  -0.93%  '),
    @('E12', '  +0.38%  '),
    @('E13', '  -3.38%  '),
    @('E14', '  +22.23%  '),
    @('E15', '  -1.37%  '),
    @('E16', '  -2.25%  '),
    @('D17', '2.567.35'),
    @('E17', '  -0.46%  '),
    @('D18', '2.284.99'),
    @('E18', '  +1.51%  '),
    @('D19', '42.400.99'),
    @('E19', '  -1.33%  '),
    @('E20', '  +6.68%  '),
    @('E21', '  -1.48%  '),
    @('E22', '  +2.09%  '),
    @('E23', '  +11.48%  '),
    @('E24', '  +0.94%  '),
    @('E25', '  +3.68%  '),
    @('E26', '  -4.60%  '),
    @('E27', '  -1.08%  '),
    @('E28', '  -5.85%  '),
    @('E29', '  -1.98%  '),
    @('E30', '  +0.59%  '),
    @('E31', '  -8.20%  '),
    @('B32', 'WEMIXToken'),
    @('C32', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('E32', '  -4.63%  '),
    @('B33', 'EthereumClassic'),
    @('C33', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('E33', '  +2.01%  '),
    @('E34', '  -3.26%  '),
    @('E35', '  +1.52%  '),
    @('E36', '  +5.71%  '),
    @('E37', '  -0.32%  '),
    @('E38', '  -3.41%  '),
    @('E39', '  -0.38%  '),
    @('E40', '  -2.76%  '),
    @('E41', '  -6.81%  '),
    @('E42', '  -0.25%  '),
    @('E43', '  -1.95%  '),
    @('E44', '  -0.30%  '),
    @('E45', '  -6.31%  '),
    @('E46', '  -2.76%  '),
    @('E47', '  -4.24%  '),
    @('E48', '  +1.40%  '),
    @('B49', 'Aave'),
    @('C49', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('E49', '  +1.96%  '),
    @('B50', 'FraxShare'),
    @('C50', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('E50', '  +0.70%  '),
    @('E51', '  -1.15%  '),
)

# Cells whose new values look numeric (e.g. "0.0920", "7.20") but must stay TEXT,
# matching the source workbook where Price/Volume columns are inline strings, not
# numbers. A leading apostrophe forces Excel to store the literal text instead of
# parsing/truncating it as a number (e.g. "0.0920" -> 0.092).
$numericLookingUpdates = @(
    @('D5', '111.82'),
    @('D6', '292.98'),
    @('D7', '0.623'),
    @('D10', '43.38'),
    @('D11', '0.0920'),
    @('D12', '54.23'),
    @('D13', '8.79'),
    @('D20', '7.20'),
    @('D22', '73.54'),
    @('D23', '3.34'),
    @('D25', '239.59'),
    @('D26', '8.87'),
    @('D28', '11.39'),
    @('D30', '175.07'),
    @('D31', '37.09'),
    @('D32', '3.12'),
    @('D33', '21.54'),
    @('D34', '0.0874'),
    @('D36', '4.97'),
    @('D38', '4.19'),
    @('D39', '0.0371'),
    @('D40', '0.103'),
    @('D41', '2.38'),
    @('D42', '70.85'),
    @('D43', '0.227'),
    @('D45', '12.31'),
    @('D49', '102.13'),
    @('D50', '8.49'),
    @('D51', '0.0975'),
)

foreach ($pair in $textUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

foreach ($pair in $numericLookingUpdates) {
    $ws.Range($pair[0]).Value = "'" + $pair[1]
}

Write-Host "Applied $($textUpdates.Count + $numericLookingUpdates.Count) cell updates"